# CU Consultar grupos y rentas
# Se termina el CU "Consultar grupos y rentas": se actualiza su estado a
# "Hecho" y su esfuerzo; se añade la opción para consultar todas las rentas
# (CU "Generar reporte de ingresos y egresos" pasa a "Hecho" con su esfuerzo
# actualizado); se actualiza la selección activa de la hoja.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# CU- 15 "Consultar grupos y rentas" (fila 19): estado "en proceso" -> "Hecho",
# esfuerzo 4 -> 8 horas.
$ws.Range("E19").Value = "Hecho"
$ws.Range("F19").Value = 8

# CU- 18 "Generar reporte de ingresos y egresos" (fila 22): estado
# "planificado" -> "Hecho", esfuerzo 1 -> 5 horas, incremento 0 -> 90.
$ws.Range("E22").Value = "Hecho"
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 90

# La celda activa de la hoja pasa a D10.
$ws.Range("D10").Select()
